$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set zoom level on the sheet view
$ws.Application.ActiveWindow.Zoom = 156

# New subject string used in replaced cells
$newSubject = "Surgery Seminar/Slide"

# Rows (block starts) where column C values change from "general surgery" to the new subject
$blockStarts = 17,39,61,83,105,127,149,171,193,215,237,259

foreach ($start in $blockStarts) {
    for ($i = 0; $i -lt 7; $i++) {
        $r = $start + $i
        $cell = $ws.Range("C$r")
        if ($i % 2 -eq 1) {
            # These rows previously used the gray-fill style (s=2); after the
            # text change they adopt the plain (no-fill) style used by their
            # neighboring odd rows (s=6). Copy that format over first.
            $srcRow = $start + ($i - 1)
            $ws.Range("C$srcRow").Copy()
            $cell.PasteSpecial(-4122)
        }
        $cell.Value = $newSubject
    }
}

$excel.CutCopyMode = 0
